$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.932.63"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "3.238.01"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "3.237.42"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  -2.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.391"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").Value = "3.793.69"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").Value = "65.053.35"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.68%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.235.38"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000159"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "413.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.71%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.205"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.495"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000111"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.77%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -2.91%  "
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("D39").Value = "2.829.56"
$ws.Range("E39").Value = "  +2.72%  "
$ws.Range("E40").Value = "  -1.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.731"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0627"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "306.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0263"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("E51").Value = "  -0.16%  "
